$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0" and J1 = "IF", matching the style of H1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-51: add I and J numeric columns ---
$I_vals = @(6,3,5,4,6,3,7,6,4,3,7,8,9,7,8,6,5,7,7,5,7,5,8,9,7,5,9,7,7,8,7,8,3,7,8,4,7,8,8,6,10,8,9,8,9,8,7,9,7,5)
$J_vals = @(6,4,5,4,6,3,7,7,4,3,7,8,9,8,8,6,5,7,7,5,7,5,8,9,7,6,9,7,7,8,7,8,4,7,8,4,7,8,8,6,10,8,9,8,9,8,7,9,7,5)

for ($idx = 0; $idx -lt 50; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $I_vals[$idx]
    $ws.Cells.Item($r, 10).Value = $J_vals[$idx]
}
